$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.117.49"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "3.128.50"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'579.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'177.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.127.86"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").Value = "'6.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").Value = "'36.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "3.650.42"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "67.081.58"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "'17.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").Value = "3.129.86"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "'490.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'0.696"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").Value = "'83.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "'12.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("D26").Value = "'2.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D29").Value = "'8.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").Value = "'2.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").Value = "'28.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("D34").Value = "0.0₃0947"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'48.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("D37").Value = "'5.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("D38").Value = "'0.948"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("D39").Value = "'49.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.311"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "'8.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("E44").Value = "  +2.83%  "
$ws.Range("D45").Value = "2.808.89"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "'375.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "'134.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  +1.74%  "
